$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style changes -------------------------------------------------------

# Header row (A1:B1): make font bold and drop the explicit "left" horizontal
# alignment (goes back to General).
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = 1

# Data row (A2:B2): font color goes from explicit black RGB to theme color 1
# (Text 1), and the fill goes from theme color 0 (white) to an explicit white
# RGB value.
$ws.Range("A2:B2").Font.ThemeColor = 1
$ws.Range("A2:B2").Interior.Color = 16777215
$ws.Range("A2:B2").Interior.PatternColor = 16777215

# --- New data row ----------------------------------------------------------

$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "admin12"

# Copy the (now updated) formatting of row 2 down onto the new row 3 so both
# rows share the same style.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
